$d = $word.ActiveDocument

$pairs = @(
    @("925×2=", "826×4="),
    @("890×2=", "470×7="),
    @("371×2=", "292×5="),
    @("978×7=", "657×8="),
    @("471×6=", "777×7="),
    @("461×3=", "253×9="),
    @("112×8=", "731×8="),
    @("837×2=", "825×4="),
    @("410×8=", "186×3="),
    @("453×8=", "631×7="),
    @("154×2=", "341×8="),
    @("480×8=", "586×5="),
    @("660×8=", "181×3="),
    @("566×4=", "209×6="),
    @("639×6=", "885×4="),
    @("330×2=", "115×2="),
    @("233×9=", "555×5="),
    @("178×6=", "824×4="),
    @("701×9=", "403×6="),
    @("190×6=", "651×7="),
    @("262×6=", "276×8="),
    @("767×6=", "772×7="),
    @("819×5=", "664×9="),
    @("910×9=", "323×5="),
    @("347×2=", "929×8=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
